$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old "regreso" column (E) on rows 1-6, keeping styles intact ---
$ws.Range("E1:E6").ClearContents()

# --- Adjust row heights for the "entrada" block (rows 2-6) ---
$ws.Rows.Item(2).RowHeight = 27
$ws.Rows.Item(3).RowHeight = 40.2
$ws.Rows.Item(4).RowHeight = 40.2
$ws.Rows.Item(5).RowHeight = 27
$ws.Rows.Item(6).RowHeight = 40.2

# --- Build the new "salida" block (rows 7-11) ---
$ws.Rows.Item(7).RowHeight = 53.4
$ws.Range("A7").Value = "principal"
$ws.Range("B7").Value = "salida"
$ws.Range("C7").Value = "`"si se encuentra en la salida de la habitaci$([char]0xF3)n principal, camine de frente por la pared izquierda hasta encontrar la primera puerta.`""
$ws.Range("D7").Value = "`"Cuidado, piso resbaladizo.`""

$ws.Rows.Item(8).RowHeight = 53.4
$ws.Range("A8").Value = "secundaria"
$ws.Range("B8").Value = "salida"
$ws.Range("C8").Value = "`"si se encuentra en la salida de la habitaci$([char]0xF3)n secundaria, gire a la izquierda y camine de frente por la pared derecha hasta encontrar la primera puerta.`""

$ws.Rows.Item(9).RowHeight = 53.4
$ws.Range("A9").Value = "cocina"
$ws.Range("B9").Value = "salida"
$ws.Range("C9").Value = "`"si se encuentra en la salida de la cocina, gire a la izquierda y camine de frente por la pared derecha hasta encontrar la primera puerta.`""
$ws.Range("D9").Value = "`"En el trayecto pueden haber materas sobre la pared derecha.`""

$ws.Rows.Item(10).RowHeight = 53.4
$ws.Range("A10").Value = "ba$([char]0xF1)o"
$ws.Range("B10").Value = "salida"
$ws.Range("C10").Value = "`"si se encuentra en la salida del ba$([char]0xF1)o, gire a la derecha y camine de frente por la pared derecha hasta encontrar la primera puerta.`""
$ws.Range("D10").Value = "`"En el trayecto hay una pared en zig zag.`""

$ws.Rows.Item(11).RowHeight = 40.2
$ws.Range("A11").Value = "patio"
$ws.Range("B11").Value = "salida"
$ws.Range("C11").Value = "`"si se encuentra en la salida del patio, camine de frente por la pared derecha hasta encontrar la primera puerta.`""
$ws.Range("D11").Value = "`"En el trayecto hay una pared en zig zag y cuadros en la pared.`""

# --- Update the active cell / selection shown when the sheet is opened ---
[void]$ws.Range("E8").Select()
